$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 4551.9546
$ws.Range("I15").Value = 4551.9546
$ws.Range("K15").Value = 13655.8638
$ws.Range("M15").Value = -13486.8638
# Row 51
$ws.Range("H51").Value = 2283.1667
$ws.Range("I51").Value = 1700
$ws.Range("J51").Value = 2866.3333
$ws.Range("K51").Value = 1700
$ws.Range("L51").Value = 2866.3333
$ws.Range("M51").Value = -1216
$ws.Range("N51").Value = -3834.3333
# Row 53
$ws.Range("H53").Value = 418.2903
$ws.Range("I53").Value = 487
$ws.Range("J53").Value = 61
$ws.Range("K53").Value = 487
$ws.Range("L53").Value = 61
$ws.Range("M53").Value = 150
$ws.Range("N53").Value = -1335
# Row 92
$ws.Range("H92").Value = 748.64703
$ws.Range("I92").Value = 615.13336
$ws.Range("J92").Value = 1750
$ws.Range("K92").Value = 615.13336
$ws.Range("L92").Value = 1750
$ws.Range("M92").Value = 632.86664
$ws.Range("N92").Value = -4246
# Row 98
$ws.Range("H98").Value = 10032.941
$ws.Range("I98").Value = 4704
$ws.Range("K98").Value = 4704
$ws.Range("M98").Value = -3206
# Row 107
$ws.Range("H107").Value = 2504.375
$ws.Range("I107").Value = 2158.7693
$ws.Range("J107").Value = 4002
$ws.Range("K107").Value = 2158.7693
$ws.Range("L107").Value = 4002
$ws.Range("M107").Value = -238.7692999999999
$ws.Range("N107").Value = -7842
# Row 122
$ws.Range("H122").Value = 10032.941
$ws.Range("I122").Value = 4704
$ws.Range("K122").Value = 14112
$ws.Range("M122").Value = -11662

$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 825
$ws.Range("I4").Value = 700
$ws.Range("J4").Value = 950
$ws.Range("K4").Value = 700
$ws.Range("L4").Value = 950
$ws.Range("M4").Value = -584
$ws.Range("N4").Value = -1182
# Row 5
$ws.Range("H5").Value = 185.75
$ws.Range("I5").Value = 185.75
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 185.75
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -73.75
$ws.Range("N5").Value = $null
# Row 102
$ws.Range("H102").Value = 33334254
$ws.Range("I102").Value = 41667544
$ws.Range("J102").Value = 1100
$ws.Range("K102").Value = 41667544
$ws.Range("L102").Value = 1100
$ws.Range("M102").Value = -41665922
$ws.Range("N102").Value = -4344
# Row 132
$ws.Range("H132").Value = 2244.6667
$ws.Range("I132").Value = 1886.1428
$ws.Range("K132").Value = 5658.428400000001
$ws.Range("M132").Value = -3128.428400000001

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 185.75
$ws.Range("I4").Value = 185.75
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 185.75
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -70.75
$ws.Range("N4").Value = $null
# Row 99
$ws.Range("H99").Value = 55557110
$ws.Range("I99").Value = 71430000
$ws.Range("K99").Value = 71430000
$ws.Range("M99").Value = -71428502
# Row 107
$ws.Range("H107").Value = 1744.909
$ws.Range("I107").Value = 1843
$ws.Range("J107").Value = 1573.25
$ws.Range("K107").Value = 1843
$ws.Range("L107").Value = 1573.25
$ws.Range("M107").Value = 77
$ws.Range("N107").Value = -5413.25
# Row 134
$ws.Range("H134").Value = 4357.6216
$ws.Range("I134").Value = 1106.1923
$ws.Range("K134").Value = 3318.5769
$ws.Range("M134").Value = -783.5769

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1227.5741
$ws.Range("I31").Value = 1198.2449
$ws.Range("J31").Value = 1515
$ws.Range("K31").Value = 1198.2449
$ws.Range("L31").Value = 1515
$ws.Range("M31").Value = -903.2448999999999
$ws.Range("N31").Value = -2105
# Row 34
$ws.Range("H34").Value = 1227.5741
$ws.Range("I34").Value = 1198.2449
$ws.Range("J34").Value = 1515
$ws.Range("K34").Value = 1198.2449
$ws.Range("L34").Value = 1515
$ws.Range("M34").Value = -996.2448999999999
$ws.Range("N34").Value = -1919
# Row 50
$ws.Range("H50").Value = 27000
$ws.Range("J50").Value = 27000
$ws.Range("L50").Value = 27000
$ws.Range("N50").Value = -28250
# Row 58
$ws.Range("H58").Value = 853.5517
$ws.Range("I58").Value = 792.6818
$ws.Range("K58").Value = 792.6818
$ws.Range("M58").Value = -589.6818
# Row 94
$ws.Range("H94").Value = 1001.7273
$ws.Range("J94").Value = 1123.2858
$ws.Range("L94").Value = 1123.2858
$ws.Range("N94").Value = -2025.2858
# Row 136
$ws.Range("H136").Value = 853.5517
$ws.Range("I136").Value = 792.6818
$ws.Range("K136").Value = 2378.0454
$ws.Range("M136").Value = 171.9546

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2996.875
$ws.Range("I80").Value = 2675
$ws.Range("J80").Value = 3190
$ws.Range("K80").Value = 2675
$ws.Range("L80").Value = 3190
$ws.Range("M80").Value = -1677
$ws.Range("N80").Value = -5186
# Row 83
$ws.Range("H83").Value = 2996.875
$ws.Range("I83").Value = 2675
$ws.Range("J83").Value = 3190
$ws.Range("K83").Value = 13375
$ws.Range("L83").Value = 15950
$ws.Range("M83").Value = -8383
$ws.Range("N83").Value = -25934
# Row 132
$ws.Range("H132").Value = 1939.8372
$ws.Range("I132").Value = 1400.3846
$ws.Range("J132").Value = 2764.8823
$ws.Range("K132").Value = 4201.1538
$ws.Range("L132").Value = 8294.6469
$ws.Range("M132").Value = -1671.1538
$ws.Range("N132").Value = -13354.6469

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1982.9166
$ws.Range("I7").Value = 1786.875
$ws.Range("J7").Value = 2375
$ws.Range("K7").Value = 1786.875
$ws.Range("L7").Value = 2375
$ws.Range("M7").Value = -1674.875
$ws.Range("N7").Value = -2599
# Row 126
$ws.Range("H126").Value = 1982.9166
$ws.Range("I126").Value = 1786.875
$ws.Range("J126").Value = 2375
$ws.Range("K126").Value = 5360.625
$ws.Range("L126").Value = 7125
$ws.Range("M126").Value = -2890.625
$ws.Range("N126").Value = -12065

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 1176.5682
$ws.Range("I132").Value = 1029.8788
$ws.Range("J132").Value = 1616.6364
$ws.Range("K132").Value = 3089.6364
$ws.Range("L132").Value = 4849.9092
$ws.Range("M132").Value = -9909.9092
